# Add Impediment Backlog Change Team Wiki
#
# 1) Drop the unused Sheet3.
# 2) Rename Sheet1 -> "Project Backlog", Sheet2 -> "Impediment Backlog".
# 3) Project Backlog: widen the title merge from A1:F1 to A1:G1 (G1 was
#    already implied by the header row) and refresh the view/selection.
# 4) Impediment Backlog: populate with the "Group 2 Impediment" table and
#    make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- 1) remove Sheet3 -------------------------------------------------
$wb.Worksheets.Item("Sheet3").Delete()

# --- 2) rename the remaining sheets ------------------------------------
$wsProject = $wb.Worksheets.Item("Sheet1")
$wsProject.Name = "Project Backlog"

$wsImpediment = $wb.Worksheets.Item("Sheet2")
$wsImpediment.Name = "Impediment Backlog"

# --- 3) Project Backlog: extend the title merge to include column G ----
$wsProject.Range("G1").Value = ""
$wsProject.Range("A1:G1").Merge()
$wsProject.Range("A1:G1").Select()

# --- 4) Impediment Backlog: build the table -----------------------------
$wsImpediment.Columns.Item(2).ColumnWidth = 10.3
$wsImpediment.Columns.Item(3).ColumnWidth = 44.6

$wsImpediment.Range("B2").Value = "Status"
$wsImpediment.Range("A1").Value = "Group 2 Impediment"
$wsImpediment.Range("A1:C1").HorizontalAlignment = -4108
$wsImpediment.Range("A2").Value = "No."
$wsImpediment.Range("C2").Value = "Description"

$wsImpediment.Range("A3").Value = 1
$wsImpediment.Range("C3").Value = "Some daily meeting late more than 20 minutes"
$wsImpediment.Range("B3").Value = "Not solve"
$wsImpediment.Range("B3").Font.Color = 255

$wsImpediment.Range("A4").Value = 2
$wsImpediment.Range("C4").Value = "Not follow convention when commit code to SVN"
$wsImpediment.Range("B4").Value = "Not solve"
$wsImpediment.Range("B4").Font.Color = 255

$wsImpediment.Range("A1:C1").Merge()
$wsImpediment.Range("C5").Select()

# --- activate the Impediment Backlog tab (also fixes activeTab / tabSelected)
$wsImpediment.Activate()
